$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Append a new data row (row 50) for the 2026-01-13 run, matching the
# existing Date/Portfolio/KAS/BTC profit columns (A:J).
$row = 50

# Column A holds dates stored as plain text (e.g. "01/12/2026"), not Excel
# date serials, so force text formatting before assigning the value to
# keep it from being auto-parsed into a date number.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "01/13/2026"

$ws.Cells.Item($row, 2).Value = 13182.45
$ws.Cells.Item($row, 3).Value = 0.2123713965796088
$ws.Cells.Item($row, 4).Value = 0.7876286034203912
$ws.Cells.Item($row, 5).Value = -109
$ws.Cells.Item($row, 6).Value = -17.17
$ws.Cells.Item($row, 7).Value = -20308.88
$ws.Cells.Item($row, 8).Value = -66.17
$ws.Cells.Item($row, 9).Value = -353.27
$ws.Cells.Item($row, 10).Value = -11.2
